$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 7 entered fully first (Name + Value)
$ws.Range("A7").Value = "EmailSubject"
$ws.Range("B7").Value = "Statement "

# Then the "Name" column filled down through row 12 (including the Bank_Names pair)
$ws.Range("A8").Value = "HDFC_Subject"
$ws.Range("A9").Value = "ICICI_Subject"
$ws.Range("A10").Value = "AXIS_Subject"
$ws.Range("A11").Value = "SBI_Subject"
$ws.Range("A12").Value = "Bank_Names"
$ws.Range("B12").Value = "HDFC,AXIS,ICICI,SBI"

# Finally, the remaining "Value" column cells were filled in
$ws.Range("B8").Value = "HDFC Statement"
$ws.Range("B9").Value = "ICICI Statement"
$ws.Range("B10").Value = "AXIS Statement"
$ws.Range("B11").Value = "SBI Statement"

$ws.Range("B19").Select()
